$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for column I (I0) and column J (IF)
$iValues = @(9, 7, 6, 9, 6, 9, 8, 7, 8, 7, 8, 9, 7, 8, 8, 7, 7)
$jValues = @(9, 9, 7, 9, 7, 9, 9, 8, 9, 7, 8, 9, 7, 8, 8, 7, 7)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
